# Minor tweak to one slide:
# The "CVM (virtual machine target for compiler" bullet is missing its
# closing parenthesis. Split the run into
#   "CVM (virtual machine target "  +  "for compiler)"
# and fix up the ")" that was missing before.

$p = $ppt.ActivePresentation

$oldText   = "CVM (virtual machine target for compiler"
$prefix    = "CVM (virtual machine target "
$oldSuffix = "for compiler"
$newSuffix = "for compiler)"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $fullText = $tr.Text

            if ($fullText -like "*$oldText*") {
                $charIndex = $fullText.IndexOf($oldText)
                $startPos  = $charIndex + 1   # TextRange.Characters is 1-based

                # Re-assert the first part of the run (keeps its existing
                # run properties / "dirty" state) ...
                $firstPart = $tr.Characters($startPos, $prefix.Length)
                $firstPart.Text = $prefix

                # ... and replace the remainder with the corrected text that
                # now includes the closing parenthesis. Setting the text on
                # just this sub-range splits it into its own run.
                $secondPart = $tr.Characters($startPos + $prefix.Length, $oldSuffix.Length)
                $secondPart.Text = $newSuffix
            }
        }
    }
}
